$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:N1").Value = "_"
$ws.Range("P1:AF1").Value = "_"
$ws.Range("A2:AF2").Value = "_"
$ws.Range("A3:AF3").Value = "_"
$ws.Range("A4:AF4").Value = "_"
$ws.Range("A5:D5").Value = "_"
$ws.Range("AC5:AF5").Value = "_"
$ws.Range("A6:D6").Value = "_"
$ws.Range("AC6:AF6").Value = "_"
$ws.Range("A7:D7").Value = "_"
$ws.Range("AC7:AF7").Value = "_"
$ws.Range("A8:D8").Value = "_"
$ws.Range("AC8:AF8").Value = "_"
$ws.Range("A9:D9").Value = "_"
$ws.Range("AC9:AF9").Value = "_"
$ws.Range("B10:D10").Value = "_"
$ws.Range("AD10:AF10").Value = "_"
$ws.Range("A11:D11").Value = "_"
$ws.Range("AC11:AF11").Value = "_"
$ws.Range("A12:D12").Value = "_"
$ws.Range("AC12:AF12").Value = "_"
$ws.Range("A13:D13").Value = "_"
$ws.Range("AC13:AF13").Value = "_"
$ws.Range("A14:D14").Value = "_"
$ws.Range("AC14:AF14").Value = "_"
$ws.Range("A15:D15").Value = "_"
$ws.Range("AC15:AF15").Value = "_"
$ws.Range("A16:D16").Value = "_"
$ws.Range("AC16:AF16").Value = "_"
$ws.Range("A17:D17").Value = "_"
$ws.Range("AC17:AF17").Value = "_"
$ws.Range("A18:D18").Value = "_"
$ws.Range("AC18:AF18").Value = "_"
$ws.Range("A19:N19").Value = "_"
$ws.Range("P19:AF19").Value = "_"
$ws.Range("A20:AF20").Value = "_"
$ws.Range("A21:AF21").Value = "_"
$ws.Range("A22:AF22").Value = "_"

[void]$ws.Range("P1").Select()
